$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "median price" row (row 2) entirely, shifting remaining rows up.
$ws.Rows.Item(2).Delete()
